$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new value would otherwise be
# auto-coerced to a number by Excel's type inference, so they stay
# text (matching the source data, which stores every cell as a string).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

# Apply the updated values from the crypto-price refresh.
$ws.Range("D2").Value = '68.000.00'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '2.537.42'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '591.48'
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").Value = '173.91'
$ws.Range("E6").Value = '  -0.74%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '0.524'
$ws.Range("E8").Value = '  -1.08%  '
$ws.Range("D9").Value = '2.536.90'
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("D10").Value = '0.138'
$ws.Range("E10").Value = '  -1.44%  '
$ws.Range("E11").Value = '  +1.86%  '
$ws.Range("E12").Value = '  -2.30%  '
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").Value = '26.41'
$ws.Range("E14").Value = '  -1.19%  '
$ws.Range("D15").Value = '3.017.99'
$ws.Range("E15").Value = '  +0.82%  '
$ws.Range("E16").Value = '  -0.75%  '
$ws.Range("D17").Value = '68.017.68'
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("E18").Value = '  +137.47%  '
$ws.Range("D19").Value = '2.506.39'
$ws.Range("E19").Value = '  -0.39%  '
$ws.Range("D20").Value = '11.77'
$ws.Range("E20").Value = '  +2.95%  '
$ws.Range("D21").Value = '7.99'
$ws.Range("E21").Value = '  -2.09%  '
$ws.Range("D22").Value = '370.54'
$ws.Range("E22").Value = '  +3.58%  '
$ws.Range("E23").Value = '  -1.31%  '
$ws.Range("E24").Value = '  -1.66%  '
$ws.Range("D25").Value = '71.90'
$ws.Range("E25").Value = '  +2.88%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  -4.18%  '
$ws.Range("D28").Value = '9.93'
$ws.Range("E28").Value = '  -3.40%  '
$ws.Range("D29").Value = '2.666.77'
$ws.Range("E29").Value = '  +0.55%  '
$ws.Range("E30").Value = '  -2.15%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = '538.93'
$ws.Range("E31").Value = '  -2.58%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '8.31'
$ws.Range("E32").Value = '  +0.66%  '
$ws.Range("E33").Value = '  -2.25%  '
$ws.Range("E34").Value = '  +0.33%  '
$ws.Range("E35").Value = '  -1.65%  '
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").Value = '159.79'
$ws.Range("E37").Value = '  +2.01%  '
$ws.Range("E38").Value = '  -2.53%  '
$ws.Range("D39").Value = '19.20'
$ws.Range("E39").Value = '  +2.15%  '
$ws.Range("E40").Value = '  +0.75%  '
$ws.Range("D41").Value = '0.350'
$ws.Range("E41").Value = '  -1.56%  '
$ws.Range("D42").Value = '5.14'
$ws.Range("E42").Value = '  -0.32%  '
$ws.Range("E43").Value = '  -1.68%  '
$ws.Range("E44").Value = '  -1.64%  '
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("D46").Value = '39.40'
$ws.Range("E46").Value = '  -1.11%  '
$ws.Range("D47").Value = '0.0₆0286'
$ws.Range("E47").Value = '  +2.72%  '
$ws.Range("D48").Value = '148.01'
$ws.Range("E48").Value = '  -0.97%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("E50").Value = '  -1.70%  '
$ws.Range("E51").Value = '  +0.93%  '
